$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 252, shifting all
# subsequent rows (252-306) down by one (to 253-307).
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range('A252').Value = 4
$ws.Range('B252').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C252').Value = 'Los Lagos'
$ws.Range('D252').Value = 44785
$ws.Range('E252').Value = 10
$ws.Range('F252').Value = 100112043
$ws.Range('G252').Value = 'Pepino ensalada'
$ws.Range('H252').Value = 'Sin especificar'
$ws.Range('I252').Value = 'Primera'
$ws.Range('J252').Value = 300
$ws.Range('K252').Value = 27000
$ws.Range('L252').Value = 27000
$ws.Range('M252').Value = 27000
$ws.Range('N252').Value = '$/caja 60 unidades'
$ws.Range('O252').Value = 'Región de Arica y Parinacota'
$ws.Range('P252').Value = 450
$ws.Range('Q252').Value = 60
$ws.Range('R252').Value = 'Hortaliza'
